# "include new fields in medication views"
#
# Re-layout the header row so that it exposes the fields expected by the
# medication import views: name, address, main_phone, phone, email, web,
# notes. The old header row used C1:I1 (telefono, Citas, email, web,
# Observaciones, Usuario, Password); the new one uses A1:G1 and drops the
# Usuario/Password columns entirely while adding name/address up front and
# renaming telefono/Citas -> main_phone/phone and Observaciones -> notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the now-unused Usuario/Password header columns.
$ws.Range("H1:I1").ClearContents()

# Write out the new header row across A1:G1.
$ws.Cells.Item(1, 1).Value = "name"
$ws.Cells.Item(1, 2).Value = "address"
$ws.Cells.Item(1, 3).Value = "main_phone"
$ws.Cells.Item(1, 4).Value = "phone"
$ws.Cells.Item(1, 5).Value = "email"
$ws.Cells.Item(1, 6).Value = "web"
$ws.Cells.Item(1, 7).Value = "notes"

# Tidy up stray leading/trailing whitespace left in a handful of data cells.
$ws.Range("B2").Value = "Calle Gran Via Del Este, Nº 80 28031, Madrid"
$ws.Range("C3").Value = "91 839 40 00"
$ws.Range("D3").Value = "91 839 45 18"
$ws.Range("B4").Value = "Avenida Del Ventisquero De La Condesa, Nº 42" + [char]10 + "28035, Madrid"
$ws.Range("B5").Value = "Avenida De Valladolid, Nº 83 28008, Madrid"
$ws.Range("B7").Value = "Calle De Modesto Lafuente, Nº 14 28010," + [char]10 + "Madrid"
$ws.Range("A9").Value = "Hospital" + [char]10 + "Clinico San Carlos"
$ws.Range("B9").Value = "Calle Del Profesor Martín Lagos, Nº S/N 28040," + [char]10 + "Madrid"
$ws.Range("B10").Value = "Calle De Serrano, Nº 199 28016, Madrid"
$ws.Range("B11").Value = "Paseo De La Castellana, Nº 261 28046, Madrid"
$ws.Range("B12").Value = "Calle Del Doctor Esquerdo, Nº 46 28009, Madrid"
$ws.Range("C14").Value = "916 70 02 57"

# Reflect the new used range now that columns H:I are empty, and move the
# active selection the way the source workbook has it after the edit.
$ws.Range("I1").Select()
